# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5264
$ws1.Range("F3").Value = 571
$ws1.Range("F4").Value = 10673
$ws1.Range("F6").Value = 572
$ws1.Range("F7").Value = 145
$ws1.Range("F8").Value = 170
$ws1.Range("F9").Value = 873

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5264
$ws4.Range("F5").Value = 571
$ws4.Range("F7").Value = 10673
$ws4.Range("F9").Value = 572
$ws4.Range("F10").Value = 145
$ws4.Range("F13").Value = 170
$ws4.Range("F14").Value = 873
